$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot original B,D,E,F,G,H values per row before any writes, since this change
# permutes the content of existing rows (B/D/E/F/G/H) while A and C stay fixed.
$orig = @{}
$orig[2] = @{
  B = $ws.Range("B2").Value2
  D = $ws.Range("D2").Value2
  E = $ws.Range("E2").Value2
  F = $ws.Range("F2").Value2
  G = $ws.Range("G2").Value2
  H = $ws.Range("H2").Value2
}
$orig[3] = @{
  B = $ws.Range("B3").Value2
  D = $ws.Range("D3").Value2
  E = $ws.Range("E3").Value2
  F = $ws.Range("F3").Value2
  G = $ws.Range("G3").Value2
  H = $ws.Range("H3").Value2
}
$orig[4] = @{
  B = $ws.Range("B4").Value2
  D = $ws.Range("D4").Value2
  E = $ws.Range("E4").Value2
  F = $ws.Range("F4").Value2
  G = $ws.Range("G4").Value2
  H = $ws.Range("H4").Value2
}
$orig[5] = @{
  B = $ws.Range("B5").Value2
  D = $ws.Range("D5").Value2
  E = $ws.Range("E5").Value2
  F = $ws.Range("F5").Value2
  G = $ws.Range("G5").Value2
  H = $ws.Range("H5").Value2
}
$orig[6] = @{
  B = $ws.Range("B6").Value2
  D = $ws.Range("D6").Value2
  E = $ws.Range("E6").Value2
  F = $ws.Range("F6").Value2
  G = $ws.Range("G6").Value2
  H = $ws.Range("H6").Value2
}
$orig[7] = @{
  B = $ws.Range("B7").Value2
  D = $ws.Range("D7").Value2
  E = $ws.Range("E7").Value2
  F = $ws.Range("F7").Value2
  G = $ws.Range("G7").Value2
  H = $ws.Range("H7").Value2
}
$orig[8] = @{
  B = $ws.Range("B8").Value2
  D = $ws.Range("D8").Value2
  E = $ws.Range("E8").Value2
  F = $ws.Range("F8").Value2
  G = $ws.Range("G8").Value2
  H = $ws.Range("H8").Value2
}
$orig[9] = @{
  B = $ws.Range("B9").Value2
  D = $ws.Range("D9").Value2
  E = $ws.Range("E9").Value2
  F = $ws.Range("F9").Value2
  G = $ws.Range("G9").Value2
  H = $ws.Range("H9").Value2
}
$orig[10] = @{
  B = $ws.Range("B10").Value2
  D = $ws.Range("D10").Value2
  E = $ws.Range("E10").Value2
  F = $ws.Range("F10").Value2
  G = $ws.Range("G10").Value2
  H = $ws.Range("H10").Value2
}
$orig[11] = @{
  B = $ws.Range("B11").Value2
  D = $ws.Range("D11").Value2
  E = $ws.Range("E11").Value2
  F = $ws.Range("F11").Value2
  G = $ws.Range("G11").Value2
  H = $ws.Range("H11").Value2
}
$orig[12] = @{
  B = $ws.Range("B12").Value2
  D = $ws.Range("D12").Value2
  E = $ws.Range("E12").Value2
  F = $ws.Range("F12").Value2
  G = $ws.Range("G12").Value2
  H = $ws.Range("H12").Value2
}
$orig[13] = @{
  B = $ws.Range("B13").Value2
  D = $ws.Range("D13").Value2
  E = $ws.Range("E13").Value2
  F = $ws.Range("F13").Value2
  G = $ws.Range("G13").Value2
  H = $ws.Range("H13").Value2
}
$orig[14] = @{
  B = $ws.Range("B14").Value2
  D = $ws.Range("D14").Value2
  E = $ws.Range("E14").Value2
  F = $ws.Range("F14").Value2
  G = $ws.Range("G14").Value2
  H = $ws.Range("H14").Value2
}
$orig[15] = @{
  B = $ws.Range("B15").Value2
  D = $ws.Range("D15").Value2
  E = $ws.Range("E15").Value2
  F = $ws.Range("F15").Value2
  G = $ws.Range("G15").Value2
  H = $ws.Range("H15").Value2
}
$orig[16] = @{
  B = $ws.Range("B16").Value2
  D = $ws.Range("D16").Value2
  E = $ws.Range("E16").Value2
  F = $ws.Range("F16").Value2
  G = $ws.Range("G16").Value2
  H = $ws.Range("H16").Value2
}
$orig[17] = @{
  B = $ws.Range("B17").Value2
  D = $ws.Range("D17").Value2
  E = $ws.Range("E17").Value2
  F = $ws.Range("F17").Value2
  G = $ws.Range("G17").Value2
  H = $ws.Range("H17").Value2
}
$orig[18] = @{
  B = $ws.Range("B18").Value2
  D = $ws.Range("D18").Value2
  E = $ws.Range("E18").Value2
  F = $ws.Range("F18").Value2
  G = $ws.Range("G18").Value2
  H = $ws.Range("H18").Value2
}
$orig[19] = @{
  B = $ws.Range("B19").Value2
  D = $ws.Range("D19").Value2
  E = $ws.Range("E19").Value2
  F = $ws.Range("F19").Value2
  G = $ws.Range("G19").Value2
  H = $ws.Range("H19").Value2
}
$orig[20] = @{
  B = $ws.Range("B20").Value2
  D = $ws.Range("D20").Value2
  E = $ws.Range("E20").Value2
  F = $ws.Range("F20").Value2
  G = $ws.Range("G20").Value2
  H = $ws.Range("H20").Value2
}
$orig[21] = @{
  B = $ws.Range("B21").Value2
  D = $ws.Range("D21").Value2
  E = $ws.Range("E21").Value2
  F = $ws.Range("F21").Value2
  G = $ws.Range("G21").Value2
  H = $ws.Range("H21").Value2
}
$orig[22] = @{
  B = $ws.Range("B22").Value2
  D = $ws.Range("D22").Value2
  E = $ws.Range("E22").Value2
  F = $ws.Range("F22").Value2
  G = $ws.Range("G22").Value2
  H = $ws.Range("H22").Value2
}
$orig[23] = @{
  B = $ws.Range("B23").Value2
  D = $ws.Range("D23").Value2
  E = $ws.Range("E23").Value2
  F = $ws.Range("F23").Value2
  G = $ws.Range("G23").Value2
  H = $ws.Range("H23").Value2
}
$orig[24] = @{
  B = $ws.Range("B24").Value2
  D = $ws.Range("D24").Value2
  E = $ws.Range("E24").Value2
  F = $ws.Range("F24").Value2
  G = $ws.Range("G24").Value2
  H = $ws.Range("H24").Value2
}
$orig[25] = @{
  B = $ws.Range("B25").Value2
  D = $ws.Range("D25").Value2
  E = $ws.Range("E25").Value2
  F = $ws.Range("F25").Value2
  G = $ws.Range("G25").Value2
  H = $ws.Range("H25").Value2
}
$orig[26] = @{
  B = $ws.Range("B26").Value2
  D = $ws.Range("D26").Value2
  E = $ws.Range("E26").Value2
  F = $ws.Range("F26").Value2
  G = $ws.Range("G26").Value2
  H = $ws.Range("H26").Value2
}
$orig[27] = @{
  B = $ws.Range("B27").Value2
  D = $ws.Range("D27").Value2
  E = $ws.Range("E27").Value2
  F = $ws.Range("F27").Value2
  G = $ws.Range("G27").Value2
  H = $ws.Range("H27").Value2
}
$orig[28] = @{
  B = $ws.Range("B28").Value2
  D = $ws.Range("D28").Value2
  E = $ws.Range("E28").Value2
  F = $ws.Range("F28").Value2
  G = $ws.Range("G28").Value2
  H = $ws.Range("H28").Value2
}
$orig[29] = @{
  B = $ws.Range("B29").Value2
  D = $ws.Range("D29").Value2
  E = $ws.Range("E29").Value2
  F = $ws.Range("F29").Value2
  G = $ws.Range("G29").Value2
  H = $ws.Range("H29").Value2
}
$orig[30] = @{
  B = $ws.Range("B30").Value2
  D = $ws.Range("D30").Value2
  E = $ws.Range("E30").Value2
  F = $ws.Range("F30").Value2
  G = $ws.Range("G30").Value2
  H = $ws.Range("H30").Value2
}
$orig[31] = @{
  B = $ws.Range("B31").Value2
  D = $ws.Range("D31").Value2
  E = $ws.Range("E31").Value2
  F = $ws.Range("F31").Value2
  G = $ws.Range("G31").Value2
  H = $ws.Range("H31").Value2
}
$orig[32] = @{
  B = $ws.Range("B32").Value2
  D = $ws.Range("D32").Value2
  E = $ws.Range("E32").Value2
  F = $ws.Range("F32").Value2
  G = $ws.Range("G32").Value2
  H = $ws.Range("H32").Value2
}
$orig[33] = @{
  B = $ws.Range("B33").Value2
  D = $ws.Range("D33").Value2
  E = $ws.Range("E33").Value2
  F = $ws.Range("F33").Value2
  G = $ws.Range("G33").Value2
  H = $ws.Range("H33").Value2
}
$orig[34] = @{
  B = $ws.Range("B34").Value2
  D = $ws.Range("D34").Value2
  E = $ws.Range("E34").Value2
  F = $ws.Range("F34").Value2
  G = $ws.Range("G34").Value2
  H = $ws.Range("H34").Value2
}
$orig[35] = @{
  B = $ws.Range("B35").Value2
  D = $ws.Range("D35").Value2
  E = $ws.Range("E35").Value2
  F = $ws.Range("F35").Value2
  G = $ws.Range("G35").Value2
  H = $ws.Range("H35").Value2
}
$orig[36] = @{
  B = $ws.Range("B36").Value2
  D = $ws.Range("D36").Value2
  E = $ws.Range("E36").Value2
  F = $ws.Range("F36").Value2
  G = $ws.Range("G36").Value2
  H = $ws.Range("H36").Value2
}
$orig[37] = @{
  B = $ws.Range("B37").Value2
  D = $ws.Range("D37").Value2
  E = $ws.Range("E37").Value2
  F = $ws.Range("F37").Value2
  G = $ws.Range("G37").Value2
  H = $ws.Range("H37").Value2
}
$orig[38] = @{
  B = $ws.Range("B38").Value2
  D = $ws.Range("D38").Value2
  E = $ws.Range("E38").Value2
  F = $ws.Range("F38").Value2
  G = $ws.Range("G38").Value2
  H = $ws.Range("H38").Value2
}
$orig[39] = @{
  B = $ws.Range("B39").Value2
  D = $ws.Range("D39").Value2
  E = $ws.Range("E39").Value2
  F = $ws.Range("F39").Value2
  G = $ws.Range("G39").Value2
  H = $ws.Range("H39").Value2
}
$orig[40] = @{
  B = $ws.Range("B40").Value2
  D = $ws.Range("D40").Value2
  E = $ws.Range("E40").Value2
  F = $ws.Range("F40").Value2
  G = $ws.Range("G40").Value2
  H = $ws.Range("H40").Value2
}
$orig[41] = @{
  B = $ws.Range("B41").Value2
  D = $ws.Range("D41").Value2
  E = $ws.Range("E41").Value2
  F = $ws.Range("F41").Value2
  G = $ws.Range("G41").Value2
  H = $ws.Range("H41").Value2
}

# Helper to write a value into G/H columns as TEXT (these columns hold position
# lists like "1, 5" or single digits like "3" that must stay text, not become numbers).
function Set-TextCell($ws, $ref, $val) {
  if ($null -eq $val) { $val = "" }
  $ws.Range($ref).NumberFormat = "@"
  $ws.Range($ref).Value = $val
  $ws.Range($ref).Style = "Normal"
}

# Apply the permutation: new row r gets the B/D/E/F/G/H content of old row mapping[r]
$ws.Range("B3").Value = $orig[5].B
$ws.Range("D3").Value = $orig[5].D
$ws.Range("E3").Value = $orig[5].E
$ws.Range("F3").Value = $orig[5].F
Set-TextCell $ws "G3" $orig[5].G
Set-TextCell $ws "H3" $orig[5].H
$ws.Range("B4").Value = $orig[10].B
$ws.Range("D4").Value = $orig[10].D
$ws.Range("E4").Value = $orig[10].E
$ws.Range("F4").Value = $orig[10].F
Set-TextCell $ws "G4" $orig[10].G
Set-TextCell $ws "H4" $orig[10].H
$ws.Range("B5").Value = $orig[3].B
$ws.Range("D5").Value = $orig[3].D
$ws.Range("E5").Value = $orig[3].E
$ws.Range("F5").Value = $orig[3].F
Set-TextCell $ws "G5" $orig[3].G
Set-TextCell $ws "H5" $orig[3].H
$ws.Range("B6").Value = $orig[7].B
$ws.Range("D6").Value = $orig[7].D
$ws.Range("E6").Value = $orig[7].E
$ws.Range("F6").Value = $orig[7].F
Set-TextCell $ws "G6" $orig[7].G
Set-TextCell $ws "H6" $orig[7].H
$ws.Range("B7").Value = $orig[4].B
$ws.Range("D7").Value = $orig[4].D
$ws.Range("E7").Value = $orig[4].E
$ws.Range("F7").Value = $orig[4].F
Set-TextCell $ws "G7" $orig[4].G
Set-TextCell $ws "H7" $orig[4].H
$ws.Range("B9").Value = $orig[6].B
$ws.Range("D9").Value = $orig[6].D
$ws.Range("E9").Value = $orig[6].E
$ws.Range("F9").Value = $orig[6].F
Set-TextCell $ws "G9" $orig[6].G
Set-TextCell $ws "H9" $orig[6].H
$ws.Range("B10").Value = $orig[9].B
$ws.Range("D10").Value = $orig[9].D
$ws.Range("E10").Value = $orig[9].E
$ws.Range("F10").Value = $orig[9].F
Set-TextCell $ws "G10" $orig[9].G
Set-TextCell $ws "H10" $orig[9].H
$ws.Range("B11").Value = $orig[15].B
$ws.Range("D11").Value = $orig[15].D
$ws.Range("E11").Value = $orig[15].E
$ws.Range("F11").Value = $orig[15].F
Set-TextCell $ws "G11" $orig[15].G
Set-TextCell $ws "H11" $orig[15].H
$ws.Range("B12").Value = $orig[25].B
$ws.Range("D12").Value = $orig[25].D
$ws.Range("E12").Value = $orig[25].E
$ws.Range("F12").Value = $orig[25].F
Set-TextCell $ws "G12" $orig[25].G
Set-TextCell $ws "H12" $orig[25].H
$ws.Range("B13").Value = $orig[14].B
$ws.Range("D13").Value = $orig[14].D
$ws.Range("E13").Value = $orig[14].E
$ws.Range("F13").Value = $orig[14].F
Set-TextCell $ws "G13" $orig[14].G
Set-TextCell $ws "H13" $orig[14].H
$ws.Range("B14").Value = $orig[18].B
$ws.Range("D14").Value = $orig[18].D
$ws.Range("E14").Value = $orig[18].E
$ws.Range("F14").Value = $orig[18].F
Set-TextCell $ws "G14" $orig[18].G
Set-TextCell $ws "H14" $orig[18].H
$ws.Range("B15").Value = $orig[13].B
$ws.Range("D15").Value = $orig[13].D
$ws.Range("E15").Value = $orig[13].E
$ws.Range("F15").Value = $orig[13].F
Set-TextCell $ws "G15" $orig[13].G
Set-TextCell $ws "H15" $orig[13].H
$ws.Range("B16").Value = $orig[17].B
$ws.Range("D16").Value = $orig[17].D
$ws.Range("E16").Value = $orig[17].E
$ws.Range("F16").Value = $orig[17].F
Set-TextCell $ws "G16" $orig[17].G
Set-TextCell $ws "H16" $orig[17].H
$ws.Range("B17").Value = $orig[19].B
$ws.Range("D17").Value = $orig[19].D
$ws.Range("E17").Value = $orig[19].E
$ws.Range("F17").Value = $orig[19].F
Set-TextCell $ws "G17" $orig[19].G
Set-TextCell $ws "H17" $orig[19].H
$ws.Range("B18").Value = $orig[11].B
$ws.Range("D18").Value = $orig[11].D
$ws.Range("E18").Value = $orig[11].E
$ws.Range("F18").Value = $orig[11].F
Set-TextCell $ws "G18" $orig[11].G
Set-TextCell $ws "H18" $orig[11].H
$ws.Range("B19").Value = $orig[26].B
$ws.Range("D19").Value = $orig[26].D
$ws.Range("E19").Value = $orig[26].E
$ws.Range("F19").Value = $orig[26].F
Set-TextCell $ws "G19" $orig[26].G
Set-TextCell $ws "H19" $orig[26].H
$ws.Range("B20").Value = $orig[23].B
$ws.Range("D20").Value = $orig[23].D
$ws.Range("E20").Value = $orig[23].E
$ws.Range("F20").Value = $orig[23].F
Set-TextCell $ws "G20" $orig[23].G
Set-TextCell $ws "H20" $orig[23].H
$ws.Range("B21").Value = $orig[22].B
$ws.Range("D21").Value = $orig[22].D
$ws.Range("E21").Value = $orig[22].E
$ws.Range("F21").Value = $orig[22].F
Set-TextCell $ws "G21" $orig[22].G
Set-TextCell $ws "H21" $orig[22].H
$ws.Range("B22").Value = $orig[16].B
$ws.Range("D22").Value = $orig[16].D
$ws.Range("E22").Value = $orig[16].E
$ws.Range("F22").Value = $orig[16].F
Set-TextCell $ws "G22" $orig[16].G
Set-TextCell $ws "H22" $orig[16].H
$ws.Range("B23").Value = $orig[20].B
$ws.Range("D23").Value = $orig[20].D
$ws.Range("E23").Value = $orig[20].E
$ws.Range("F23").Value = $orig[20].F
Set-TextCell $ws "G23" $orig[20].G
Set-TextCell $ws "H23" $orig[20].H
$ws.Range("B24").Value = $orig[21].B
$ws.Range("D24").Value = $orig[21].D
$ws.Range("E24").Value = $orig[21].E
$ws.Range("F24").Value = $orig[21].F
Set-TextCell $ws "G24" $orig[21].G
Set-TextCell $ws "H24" $orig[21].H
$ws.Range("B25").Value = $orig[12].B
$ws.Range("D25").Value = $orig[12].D
$ws.Range("E25").Value = $orig[12].E
$ws.Range("F25").Value = $orig[12].F
Set-TextCell $ws "G25" $orig[12].G
Set-TextCell $ws "H25" $orig[12].H
$ws.Range("B26").Value = $orig[24].B
$ws.Range("D26").Value = $orig[24].D
$ws.Range("E26").Value = $orig[24].E
$ws.Range("F26").Value = $orig[24].F
Set-TextCell $ws "G26" $orig[24].G
Set-TextCell $ws "H26" $orig[24].H
$ws.Range("B27").Value = $orig[37].B
$ws.Range("D27").Value = $orig[37].D
$ws.Range("E27").Value = $orig[37].E
$ws.Range("F27").Value = $orig[37].F
Set-TextCell $ws "G27" $orig[37].G
Set-TextCell $ws "H27" $orig[37].H
$ws.Range("B28").Value = $orig[34].B
$ws.Range("D28").Value = $orig[34].D
$ws.Range("E28").Value = $orig[34].E
$ws.Range("F28").Value = $orig[34].F
Set-TextCell $ws "G28" $orig[34].G
Set-TextCell $ws "H28" $orig[34].H
$ws.Range("B29").Value = $orig[27].B
$ws.Range("D29").Value = $orig[27].D
$ws.Range("E29").Value = $orig[27].E
$ws.Range("F29").Value = $orig[27].F
Set-TextCell $ws "G29" $orig[27].G
Set-TextCell $ws "H29" $orig[27].H
$ws.Range("B30").Value = $orig[38].B
$ws.Range("D30").Value = $orig[38].D
$ws.Range("E30").Value = $orig[38].E
$ws.Range("F30").Value = $orig[38].F
Set-TextCell $ws "G30" $orig[38].G
Set-TextCell $ws "H30" $orig[38].H
$ws.Range("B31").Value = $orig[29].B
$ws.Range("D31").Value = $orig[29].D
$ws.Range("E31").Value = $orig[29].E
$ws.Range("F31").Value = $orig[29].F
Set-TextCell $ws "G31" $orig[29].G
Set-TextCell $ws "H31" $orig[29].H
$ws.Range("B32").Value = $orig[33].B
$ws.Range("D32").Value = $orig[33].D
$ws.Range("E32").Value = $orig[33].E
$ws.Range("F32").Value = $orig[33].F
Set-TextCell $ws "G32" $orig[33].G
Set-TextCell $ws "H32" $orig[33].H
$ws.Range("B33").Value = $orig[35].B
$ws.Range("D33").Value = $orig[35].D
$ws.Range("E33").Value = $orig[35].E
$ws.Range("F33").Value = $orig[35].F
Set-TextCell $ws "G33" $orig[35].G
Set-TextCell $ws "H33" $orig[35].H
$ws.Range("B34").Value = $orig[30].B
$ws.Range("D34").Value = $orig[30].D
$ws.Range("E34").Value = $orig[30].E
$ws.Range("F34").Value = $orig[30].F
Set-TextCell $ws "G34" $orig[30].G
Set-TextCell $ws "H34" $orig[30].H
$ws.Range("B35").Value = $orig[32].B
$ws.Range("D35").Value = $orig[32].D
$ws.Range("E35").Value = $orig[32].E
$ws.Range("F35").Value = $orig[32].F
Set-TextCell $ws "G35" $orig[32].G
Set-TextCell $ws "H35" $orig[32].H
$ws.Range("B37").Value = $orig[31].B
$ws.Range("D37").Value = $orig[31].D
$ws.Range("E37").Value = $orig[31].E
$ws.Range("F37").Value = $orig[31].F
Set-TextCell $ws "G37" $orig[31].G
Set-TextCell $ws "H37" $orig[31].H
$ws.Range("B38").Value = $orig[28].B
$ws.Range("D38").Value = $orig[28].D
$ws.Range("E38").Value = $orig[28].E
$ws.Range("F38").Value = $orig[28].F
Set-TextCell $ws "G38" $orig[28].G
Set-TextCell $ws "H38" $orig[28].H
$ws.Range("B39").Value = $orig[40].B
$ws.Range("D39").Value = $orig[40].D
$ws.Range("E39").Value = $orig[40].E
$ws.Range("F39").Value = $orig[40].F
Set-TextCell $ws "G39" $orig[40].G
Set-TextCell $ws "H39" $orig[40].H
$ws.Range("B40").Value = $orig[39].B
$ws.Range("D40").Value = $orig[39].D
$ws.Range("E40").Value = $orig[39].E
$ws.Range("F40").Value = $orig[39].F
Set-TextCell $ws "G40" $orig[39].G
Set-TextCell $ws "H40" $orig[39].H
